# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit calculation columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 376
$ws.Range("I28").Value = 162
$ws.Range("J28").Value = 732.6667
$ws.Range("K28").Value = 162
$ws.Range("L28").Value = 732.6667
$ws.Range("M28").Value = 323
$ws.Range("N28").Value = -1702.6667

$ws.Range("H43").Value = 1437.0667
$ws.Range("J43").Value = 1396.875
$ws.Range("L43").Value = 1396.875
$ws.Range("N43").Value = -1534.875

$ws.Range("H88").Value = 2619.7812
$ws.Range("J88").Value = 2738.5
$ws.Range("L88").Value = 2738.5
$ws.Range("N88").Value = -3550.5

$ws.Range("H91").Value = 2619.7812
$ws.Range("J91").Value = 2738.5
$ws.Range("L91").Value = 2738.5
$ws.Range("N91").Value = -5546.5

$ws.Range("H116").Value = 6155
$ws.Range("I116").Value = 4995
$ws.Range("J116").Value = 6445
$ws.Range("K116").Value = 4995
$ws.Range("L116").Value = 6445
$ws.Range("M116").Value = -1553
$ws.Range("N116").Value = -13329

$ws.Range("H129").Value = 50001410
$ws.Range("I129").Value = 55556784
$ws.Range("J129").Value = 3000
$ws.Range("K129").Value = 166670352
$ws.Range("L129").Value = 9000
$ws.Range("M129").Value = -166665352
$ws.Range("N129").Value = -19000

$ws.Range("H138").Value = 3007.7473
$ws.Range("J138").Value = 3150.236
$ws.Range("L138").Value = 9450.707999999999
$ws.Range("N138").Value = -19730.708

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7995219.5
$ws.Range("I45").Value = 11989601
$ws.Range("J45").Value = 6457
$ws.Range("K45").Value = 11989601
$ws.Range("L45").Value = 6457
$ws.Range("M45").Value = -11989224
$ws.Range("N45").Value = -7211

$ws.Range("H61").Value = 2321.739
$ws.Range("I61").Value = 1698.8125
$ws.Range("J61").Value = 3745.5715
$ws.Range("K61").Value = 1698.8125
$ws.Range("L61").Value = 3745.5715
$ws.Range("M61").Value = -1486.8125
$ws.Range("N61").Value = -4169.5715

$ws.Range("H112").Value = 11152.5
$ws.Range("J112").Value = 11152.5
$ws.Range("L112").Value = 11152.5
$ws.Range("N112").Value = -14106.5

$ws.Range("H124").Value = 12693.75
$ws.Range("J124").Value = 12693.75
$ws.Range("L124").Value = 12693.75
$ws.Range("N124").Value = -22513.75

$ws.Range("H136").Value = 2321.739
$ws.Range("I136").Value = 1698.8125
$ws.Range("J136").Value = 3745.5715
$ws.Range("K136").Value = 5096.4375
$ws.Range("L136").Value = 11236.7145
$ws.Range("M136").Value = -2546.4375
$ws.Range("N136").Value = -16336.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 10000
$ws.Range("J33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("N33").Value = -10672

$ws.Range("H107").Value = 2551861
$ws.Range("I107").Value = 2977004
$ws.Range("J107").Value = 1004
$ws.Range("K107").Value = 2977004
$ws.Range("L107").Value = 1004
$ws.Range("M107").Value = -2975084
$ws.Range("N107").Value = -4844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20538.25
$ws.Range("I31").Value = 1189.5312
$ws.Range("K31").Value = 1189.5312
$ws.Range("M31").Value = -894.5311999999999

$ws.Range("H34").Value = 20538.25
$ws.Range("I34").Value = 1189.5312
$ws.Range("K34").Value = 1189.5312
$ws.Range("M34").Value = -987.5311999999999

$ws.Range("H41").Value = 10779.833
$ws.Range("I41").Value = 1175
$ws.Range("K41").Value = 1175
$ws.Range("M41").Value = -747

$ws.Range("H94").Value = 987.7059
$ws.Range("I94").Value = 652.8333
$ws.Range("J94").Value = 1170.3636
$ws.Range("K94").Value = 652.8333
$ws.Range("L94").Value = 1170.3636
$ws.Range("M94").Value = -201.8333
$ws.Range("N94").Value = -2072.3636

$ws.Range("H102").Value = 51246.75
$ws.Range("J102").Value = 51246.75
$ws.Range("L102").Value = 51246.75
$ws.Range("N102").Value = -56114.75

$ws.Range("H132").Value = 101987.42
$ws.Range("I132").Value = 73473.42999999999
$ws.Range("J132").Value = 181826.6
$ws.Range("K132").Value = 220420.29
$ws.Range("L132").Value = 545479.8
$ws.Range("M132").Value = -217890.29
$ws.Range("N132").Value = -550539.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 555
$ws.Range("I3").Value = 555
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1665
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1553
$ws.Range("N3").Value = ""

$ws.Range("H63").Value = 12977.571
$ws.Range("I63").Value = 16032.75
$ws.Range("J63").Value = 8904
$ws.Range("K63").Value = 48098.25
$ws.Range("L63").Value = 26712
$ws.Range("M63").Value = -47349.25
$ws.Range("N63").Value = -28210

$ws.Range("H66").Value = 12977.571
$ws.Range("I66").Value = 16032.75
$ws.Range("J66").Value = 8904
$ws.Range("K66").Value = 144294.75
$ws.Range("L66").Value = 80136
$ws.Range("M66").Value = -140550.75
$ws.Range("N66").Value = -87624

$ws.Range("H134").Value = 2138.0667
$ws.Range("I134").Value = 2138.0667
$ws.Range("K134").Value = 6414.2001
$ws.Range("M134").Value = -1344.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 21658
$ws.Range("J15").Value = 21658
$ws.Range("L15").Value = 21658
$ws.Range("N15").Value = -22234

$ws.Range("H51").Value = 84911.11
$ws.Range("J51").Value = 84911.11
$ws.Range("L51").Value = 84911.11
$ws.Range("N51").Value = -85929.11

$ws.Range("H81").Value = 21658
$ws.Range("J81").Value = 21658
$ws.Range("L81").Value = 21658
$ws.Range("N81").Value = -23654

$ws.Range("H84").Value = 21658
$ws.Range("J84").Value = 21658
$ws.Range("L84").Value = 64974
$ws.Range("N84").Value = -74958

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6786.143
$ws.Range("I7").Value = 3252
$ws.Range("J7").Value = 8199.799999999999
$ws.Range("K7").Value = 3252
$ws.Range("L7").Value = 8199.799999999999
$ws.Range("M7").Value = -3140
$ws.Range("N7").Value = -8423.799999999999

$ws.Range("H17").Value = 25000
$ws.Range("I17").Value = 25000
$ws.Range("K17").Value = 25000
$ws.Range("M17").Value = -24830

$ws.Range("H46").Value = 4676.294
$ws.Range("I46").Value = 2980
$ws.Range("J46").Value = 5383.0835
$ws.Range("K46").Value = 2980
$ws.Range("L46").Value = 5383.0835
$ws.Range("M46").Value = -2792
$ws.Range("N46").Value = -5759.0835

$ws.Range("H110").Value = 19000
$ws.Range("J110").Value = 19000
$ws.Range("L110").Value = 19000
$ws.Range("N110").Value = -27180

$ws.Range("H126").Value = 6786.143
$ws.Range("I126").Value = 3252
$ws.Range("J126").Value = 8199.799999999999
$ws.Range("K126").Value = 9756
$ws.Range("L126").Value = 24599.4
$ws.Range("M126").Value = -7286
$ws.Range("N126").Value = -29539.4

$ws.Range("H127").Value = 64999
$ws.Range("J127").Value = 64999
$ws.Range("L127").Value = 64999
$ws.Range("N127").Value = -74919

$ws.Range("H136").Value = 62672.234
$ws.Range("I136").Value = 72132.97
$ws.Range("J136").Value = 7800
$ws.Range("K136").Value = 216398.91
$ws.Range("L136").Value = 23400
$ws.Range("M136").Value = -213848.91
$ws.Range("N136").Value = -28500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 34499.5
$ws.Range("I42").Value = 33999
$ws.Range("K42").Value = 33999
$ws.Range("M42").Value = -33621

$ws.Range("H51").Value = 13000
$ws.Range("I51").Value = 13000
$ws.Range("K51").Value = 13000
$ws.Range("M51").Value = -12490

$ws.Range("H52").Value = 11600
$ws.Range("J52").Value = 16000
$ws.Range("L52").Value = 16000
$ws.Range("N52").Value = -16452

$ws.Range("H126").Value = 2488.7334
$ws.Range("I126").Value = 3254.5
$ws.Range("J126").Value = 1613.5714
$ws.Range("K126").Value = 9763.5
$ws.Range("L126").Value = 4840.7142
$ws.Range("M126").Value = -7293.5
$ws.Range("N126").Value = -9780.7142

$ws.Range("H132").Value = 112325000
$ws.Range("I132").Value = 200002620
$ws.Range("J132").Value = 2727971.5
$ws.Range("K132").Value = 600007860
$ws.Range("L132").Value = 8183914.5
$ws.Range("M132").Value = -600005330
$ws.Range("N132").Value = -8188974.5

$ws.Range("H136").Value = 1937.8
$ws.Range("I136").Value = 1204.5862
$ws.Range("J136").Value = 5481.6665
$ws.Range("K136").Value = 3613.7586
$ws.Range("L136").Value = 16444.9995
$ws.Range("M136").Value = -1063.7586
$ws.Range("N136").Value = -21544.9995
